$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$dates = @(
    "Mon Jul 17 21:14:07 EDT 2023",
    "Mon Jul 17 21:14:15 EDT 2023",
    "Mon Jul 17 21:14:24 EDT 2023",
    "Mon Jul 17 21:14:33 EDT 2023",
    "Mon Jul 17 21:14:42 EDT 2023",
    "Mon Jul 17 21:14:52 EDT 2023",
    "Mon Jul 17 21:15:00 EDT 2023",
    "Mon Jul 17 21:15:09 EDT 2023",
    "Mon Jul 17 21:15:18 EDT 2023",
    "Mon Jul 17 21:15:27 EDT 2023",
    "Mon Jul 17 21:15:36 EDT 2023",
    "Mon Jul 17 21:15:45 EDT 2023",
    "Mon Jul 17 21:15:55 EDT 2023",
    "Mon Jul 17 21:16:04 EDT 2023",
    "Mon Jul 17 21:16:13 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
